$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 757; this shifts existing rows 757-830 down to 758-831
# and Excel automatically extends the sheet dimension (A1:T830 -> A1:T831).
$ws.Rows.Item(757).Insert()

# Populate the newly inserted row 757 with the new weekly price record.
# Columns A,B,C,E,F,G,H,I,J,K hold the same constant values used throughout the sheet.
$ws.Range("A757").Value = 11
$ws.Range("B757").Value = "Vega Monumental Concepción"
$ws.Range("C757").Value = "Bíobío"
$ws.Range("D757").Value = 45132
$ws.Range("E757").Value = 8
$ws.Range("F757").Value = "Fruta"
$ws.Range("G757").Value = 100102
$ws.Range("H757").Value = "Cítricos"
$ws.Range("I757").Value = 100102003
$ws.Range("J757").Value = "Limón"
$ws.Range("K757").Value = "Sin especificar"
$ws.Range("L757").Value = "1a amarillo"
$ws.Range("M757").Value = 270
$ws.Range("N757").Value = 6000
$ws.Range("O757").Value = 6500
$ws.Range("P757").Value = 6278
$ws.Range("Q757").Value = "$/malla 16 kilos"
$ws.Range("R757").Value = "Región de O'Higgins"
$ws.Range("S757").Value = 392
$ws.Range("T757").Value = 16
